$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AgeGroup Child -> Adult; Protocol Stimulus-response -> Single pulse
$ws.Range("E2").Value = "Adult"
$ws.Range("H2").Value = "Single pulse"

# Row 4: SequenceType rTMS -> repetitive
$ws.Range("G4").Value = "repetitive"

# Rows 8,9,10,11,12: SequenceType single_pulse -> single
$ws.Range("G8").Value = "single"
$ws.Range("G9").Value = "single"
$ws.Range("G10").Value = "single"
$ws.Range("G11").Value = "single"
$ws.Range("G12").Value = "single"

# Rows 10,11,12: Protocol Stimulus-response -> RC
$ws.Range("H10").Value = "RC"
$ws.Range("H11").Value = "RC"
$ws.Range("H12").Value = "RC"

# Row 9: remove the stray empty styled cell in column B entirely
$ws.Range("B9").Clear()

# Update the active selection to match the author's last position
$ws.Range("H12").Select()

$wb.Save()
